$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates ---------------------------------------------------
# Header row: B1 keeps "DateTime" (unchanged)
# Row 2: sample ID unchanged, but the DateTime sample switches to ISO-ish format
$ws.Range("B2").Value = "2015-06-21 17:12:10"

# Row 3: the sample employee id becomes a text value (was a pure number before)
$ws.Range("A3").Value = "14000031"
$ws.Range("B3").Value = "2015-06-04 08:12:21"

# --- Formatting updates ----------------------------------------------------
# A2 no longer wraps text
$ws.Range("A2").WrapText = $false

# Column widths
$ws.Columns.Item(2).ColumnWidth = 14.6

# Selection moves to B1
$ws.Range("B1").Select()
